# Insert a new data row at row 31 (weekly Haba price report update),
# pushing the existing rows 31..131 down to 32..132, then populate the
# new row 31 with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31..131 down by one to make room for the new observation.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new record's values.
$ws.Cells.Item(31, 1).Value = 9
$ws.Cells.Item(31, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(31, 3).Value = "Metropolitana"
$ws.Cells.Item(31, 4).Value = 44459
$ws.Cells.Item(31, 5).Value = 13
$ws.Cells.Item(31, 6).Value = 100112026
$ws.Cells.Item(31, 7).Value = "Haba"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 43
$ws.Cells.Item(31, 11).Value = 16000
$ws.Cells.Item(31, 12).Value = 17000
$ws.Cells.Item(31, 13).Value = 16488
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 660
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
